# Inserts a new data row at row 55 (pushing the existing rows 55..114 down
# to 56..115) on the active sheet, and populates the new row with the
# Maracuyá / Vega Modelo de Temuco record that the commit added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 55 and below down by one row.
$ws.Range("A55").EntireRow.Insert()

# Populate the newly inserted row 55 with the new record's values.
$ws.Range("A55").Value = 10
$ws.Range("B55").Value = "Vega Modelo de Temuco"
$ws.Range("C55").Value = "La Araucanía"
$ws.Range("D55").Value = 45159
$ws.Range("E55").Value = 9
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100108
$ws.Range("H55").Value = "Tropicales y subtropicales"
$ws.Range("I55").Value = 100108003
$ws.Range("J55").Value = "Maracuyá"
$ws.Range("K55").Value = "Sin especificar"
$ws.Range("L55").Value = "Primera"
$ws.Range("M55").Value = 80
$ws.Range("N55").Value = 38000
$ws.Range("O55").Value = 38000
$ws.Range("P55").Value = 38000
$ws.Range("Q55").Value = "$/caja 18 kilos"
$ws.Range("R55").Value = "Región de Arica y Parinacota"
$ws.Range("S55").Value = 2111
$ws.Range("T55").Value = 18
